$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1195.909
$ws.Range("J17").Value = 1140.4286
$ws.Range("L17").Value = 3421.2858
$ws.Range("N17").Value = -3757.2858
# Row 32
$ws.Range("H32").Value = 16671417
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 20004900
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 20004900
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -20005552
# Row 62
$ws.Range("H62").Value = 1959.6
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
# Row 65
$ws.Range("H65").Value = 1959.6
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
# Row 100
$ws.Range("H100").Value = 4474.5
$ws.Range("I100").Value = 3106.6365
$ws.Range("J100").Value = 5631.923
$ws.Range("K100").Value = 3106.6365
$ws.Range("L100").Value = 5631.923
$ws.Range("M100").Value = -2565.6365
$ws.Range("N100").Value = -6713.923
# Row 104
$ws.Range("H104").Value = 95.166664
$ws.Range("I104").Value = 95.166664
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 285.499992
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = 1461.500008
# Row 135
$ws.Range("H135").Value = 16668427
$ws.Range("I135").Value = 20835112
$ws.Range("K135").Value = 187516008
$ws.Range("M135").Value = -187513473
# Row 137
$ws.Range("H137").Value = 2753.0833
$ws.Range("I137").Value = 2594.2727
$ws.Range("K137").Value = 7782.8181
$ws.Range("M137").Value = -5232.8181

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15388800
$ws.Range("I32").Value = 15877300
$ws.Range("K32").Value = 15877300
$ws.Range("M32").Value = -15877013
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0
# Row 61
$ws.Range("H61").Value = 3171.889
$ws.Range("I61").Value = 3193.375
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 3193.375
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2981.375
$ws.Range("N61").Value = -3424
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0
# Row 97
$ws.Range("H97").Value = 716.7273
$ws.Range("I97").Value = 683.41174
$ws.Range("K97").Value = 683.41174
$ws.Range("M97").Value = -187.41174
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
# Row 122
$ws.Range("H122").Value = 2439.25
$ws.Range("I122").Value = 2439.25
$ws.Range("K122").Value = 7317.75
$ws.Range("M122").Value = -4867.75
# Row 132
$ws.Range("H132").Value = 2401.3914
$ws.Range("I132").Value = 2419.6365
$ws.Range("K132").Value = 7258.9095
$ws.Range("M132").Value = -4728.9095
# Row 136
$ws.Range("H136").Value = 3171.889
$ws.Range("I136").Value = 3193.375
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 9580.125
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -7030.125
$ws.Range("N136").Value = -14100
# Row 139
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2252.0356
$ws.Range("I86").Value = 2027.2222
$ws.Range("J86").Value = 2656.7
$ws.Range("K86").Value = 2027.2222
$ws.Range("L86").Value = 2656.7
$ws.Range("M86").Value = -904.2221999999999
$ws.Range("N86").Value = -4902.7
# Row 89
$ws.Range("H89").Value = 2252.0356
$ws.Range("I89").Value = 2027.2222
$ws.Range("J89").Value = 2656.7
$ws.Range("K89").Value = 10136.111
$ws.Range("L89").Value = 13283.5
$ws.Range("M89").Value = -4520.110999999999
$ws.Range("N89").Value = -24515.5
# Row 99
$ws.Range("H99").Value = 47884.445
$ws.Range("I99").Value = 53495
$ws.Range("K99").Value = 53495
$ws.Range("M99").Value = -51997
# Row 126
$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880
# Row 134
$ws.Range("H134").Value = 1791.5483
$ws.Range("I134").Value = 1617.9333
$ws.Range("K134").Value = 4853.7999
$ws.Range("M134").Value = -2318.7999
# Row 140
$ws.Range("H140").Value = 77166.414
$ws.Range("J140").Value = 77166.414
$ws.Range("L140").Value = 77166.414
$ws.Range("N140").Value = -87526.414

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 8194.666999999999
$ws.Range("I16").Value = 8494.286
$ws.Range("K16").Value = 8494.286
$ws.Range("M16").Value = -8207.286
# Row 31
$ws.Range("H31").Value = 1706.3043
$ws.Range("I31").Value = 1601.2368
$ws.Range("K31").Value = 1601.2368
$ws.Range("M31").Value = -1306.2368
# Row 34
$ws.Range("H34").Value = 1706.3043
$ws.Range("I34").Value = 1601.2368
$ws.Range("K34").Value = 1601.2368
$ws.Range("M34").Value = -1399.2368
# Row 107
$ws.Range("H107").Value = 13785.8125
$ws.Range("I107").Value = 1303
$ws.Range("J107").Value = 21275.5
$ws.Range("K107").Value = 1303
$ws.Range("L107").Value = 21275.5
$ws.Range("M107").Value = 617
$ws.Range("N107").Value = -25115.5
# Row 113
$ws.Range("H113").Value = 8194.666999999999
$ws.Range("I113").Value = 8494.286
$ws.Range("K113").Value = 8494.286
$ws.Range("M113").Value = -6324.286
# Row 132
$ws.Range("H132").Value = 2330.0588
$ws.Range("I132").Value = 2312.3125
$ws.Range("K132").Value = 6936.9375
$ws.Range("M132").Value = -4406.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 2250
$ws.Range("I107").Value = 375
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 1125
$ws.Range("L107").Value = 10500
$ws.Range("M107").Value = 795
$ws.Range("N107").Value = -14340

$ws = $wb.Worksheets.Item("GSM")
# Row 86
$ws.Range("H86").Value = 119999
$ws.Range("J86").Value = 119999
$ws.Range("L86").Value = 119999
$ws.Range("N86").Value = -122371
# Row 89
$ws.Range("H89").Value = 119999
$ws.Range("J89").Value = 119999
$ws.Range("L89").Value = 359997
$ws.Range("N89").Value = -371853
# Row 92
$ws.Range("H92").Value = 16140.143
$ws.Range("J92").Value = 17458.834
$ws.Range("L92").Value = 17458.834
$ws.Range("N92").Value = -21202.834
# Row 97
$ws.Range("H97").Value = 324.9565
$ws.Range("I97").Value = 290.05554
$ws.Range("K97").Value = 290.05554
$ws.Range("M97").Value = 205.94446

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3087.625
$ws.Range("I68").Value = 3041
$ws.Range("J68").Value = 3165.3333
$ws.Range("K68").Value = 3041
$ws.Range("L68").Value = 3165.3333
$ws.Range("M68").Value = -2292
$ws.Range("N68").Value = -4663.3333
# Row 71
$ws.Range("H71").Value = 3087.625
$ws.Range("I71").Value = 3041
$ws.Range("J71").Value = 3165.3333
$ws.Range("K71").Value = 15205
$ws.Range("L71").Value = 15826.6665
$ws.Range("M71").Value = -11461
$ws.Range("N71").Value = -23314.6665
# Row 93
$ws.Range("H93").Value = 8355.177
$ws.Range("I93").Value = 8228.4
$ws.Range("J93").Value = 8536.286
$ws.Range("K93").Value = 8228.4
$ws.Range("L93").Value = 8536.286
$ws.Range("M93").Value = -6980.4
$ws.Range("N93").Value = -11032.286
# Row 100
$ws.Range("H100").Value = 75188
$ws.Range("I100").Value = 142375.5
$ws.Range("K100").Value = 142375.5
$ws.Range("M100").Value = -141834.5
# Row 127
$ws.Range("H127").Value = 69987.09
$ws.Range("J127").Value = 69987.09
$ws.Range("L127").Value = 69987.09
$ws.Range("N127").Value = -79907.09
# Row 136
$ws.Range("H136").Value = 2982.0833
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 42202.43
$ws.Range("J69").Value = 42202.43
$ws.Range("L69").Value = 42202.43
$ws.Range("N69").Value = -43700.43
# Row 72
$ws.Range("H72").Value = 42202.43
$ws.Range("J72").Value = 42202.43
$ws.Range("L72").Value = 126607.29
$ws.Range("N72").Value = -134095.29
# Row 132
$ws.Range("H132").Value = 1694.65
$ws.Range("I132").Value = 1605.2222
$ws.Range("K132").Value = 4815.6666
$ws.Range("M132").Value = -2285.6666
# Row 136
$ws.Range("H136").Value = 1462.4
$ws.Range("I136").Value = 835.12
$ws.Range("K136").Value = 2505.36
$ws.Range("M136").Value = 44.63999999999987
